# Adds a "WebSocket Command" entry to the Commands sheet, following the
# same row layout used by the other command groups (header row with
# command name / syntax / description, followed by a spacer row that
# reuses the "sleep(<float>)" entry like the other groups do).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

$lq = [char]0x201C   # “
$rq = [char]0x201D   # ”

$ws.Range("A89").Value = "WebSocket Command"
$ws.Range("B89").Value = "send(<json>)"
$ws.Range("C89").Value = "If {} substitutions are used, json brackets need to be duplicated to escape them like in send({{ ${lq}value${rq}: {}}})"

$ws.Range("B90").Value = "sleep(<float>)"
$ws.Range("C90").Value = "sleep: add a delay of <float> seconds"

$ws.Range("C89").Select() | Out-Null
